# Radjeno na proveri unosa podataka na Preferences stranici za Urgent days
#
# Updates the "AdministrationPreferencesCertificateAtPoint" test row:
#  - D15 now lists both the FALSE and TRUE variant test names (wrapped text)
#  - E15 / E16 are marked as DONE
#  - Column D is widened to fit the new two-line content
#  - Selection / scroll position moved further down the sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# D15: replace the single test name with the FALSE/TRUE pair on two lines,
# and turn wrap-text on so both lines are visible.
$ws.Range("D15").Value = "AdministrationPreferencesCertificateAtPoint-FALSE" + [char]10 + "AdministrationPreferencesCertificateAtPoint-TRUE"
$ws.Range("D15").WrapText = $true

# Mark the CertificateAtPoint (D15) and RedactionStatement (D16) rows as DONE.
$ws.Range("E15").Value = "DONE"
$ws.Range("E16").Value = "DONE"

# Widen column D now that it holds two lines of text instead of one.
$ws.Columns.Item(4).ColumnWidth = 50.166666666666664

# Move the active selection down to where the work continued.
$ws.Activate()
$ws.Range("E17").Select()
